# Applies the "Updated symbol list" commit: refreshed coin prices and
# restored the CEJI / BKEXToken row ordering with their current values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextValue "D2" "269.39"
Set-TextValue "D3" "22.94"
Set-TextValue "D4" "6.332"
Set-TextValue "D5" "0.06195"
Set-TextValue "D6" "3.642"
Set-TextValue "D7" "6.687"
Set-TextValue "D8" "1.386"
Set-TextValue "D9" "0.8298"
Set-TextValue "D10" "0.01378"
Set-TextValue "D12" "0.08270"
Set-TextValue "D13" "0.03501"
Set-TextValue "D14" "0.03177"
Set-TextValue "D15" "0.09348"
Set-TextValue "D16" "3.840"
Set-TextValue "D17" "0.001642"
Set-TextValue "D18" "0.04743"
Set-TextValue "D19" "0.006371"
Set-TextValue "D20" "0.005677"
Set-TextValue "D21" "0.001077"
Set-TextValue "D23" "3.717"
Set-TextValue "D27" "0.0002703"
Set-TextValue "D41" "0.006981"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1164"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003320"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
Set-TextValue "D44" "0.01193"
Set-TextValue "D45" "0.00006246"
Set-TextValue "D46" "0.0009897"
Set-TextValue "D48" "0.9198"
Set-TextValue "D49" "0.002227"
